$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8
$ws.Range("AK8").Value = "te"

# Row 12
$ws.Range("V12").Value = 7001
$ws.Range("AI12").Value = "HOLD"
$ws.Range("AJ12").Value = "HOLD"
$ws.Range("AK12").Value = "testing purpose"
$ws.Range("AL12").Value = "testing purpos"
$ws.Range("AM12").Value = "testing purpose"
$ws.Range("AN12").Value = "testing purpose"
$ws.Range("AO12").Value = "testing purpose"

# Row 13
$ws.Range("AK13").Value = "te"
$ws.Range("AL13").Value = "te"
$ws.Range("AM13").Value = "te"
$ws.Range("AN13").Value = "te"
$ws.Range("AO13").Value = "te"

# Row 18
$ws.Range("AN18").Value = "te"

# Row 19
$ws.Range("AN19").Value = "te"

# Row 20
$ws.Range("AN20").Value = "te"

# Row 21
$ws.Range("AN21").Value = "te"

# Row 22
$ws.Range("AN22").Value = "te"

# Row 23
$ws.Range("AL23").Value = "testing"
$ws.Range("AM23").Value = "testing"
$ws.Range("AN23").Value = "te"

# Row 24
$ws.Range("AL24").Value = "testing"
$ws.Range("AM24").Value = "testing"
$ws.Range("AN24").Value = "te"

# Row 25
$ws.Range("AK25").Value = "te"
$ws.Range("AL25").Value = "testing"

# Row 26
$ws.Range("AK26").Value = "te"
$ws.Range("AL26").Value = "testing"

# Row 27
$ws.Range("AK27").Value = "te"
$ws.Range("AL27").Value = "te"
$ws.Range("AM27").Value = "te"
$ws.Range("AN27").Value = "te"
$ws.Range("AO27").Value = "te"
